$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 22: new match data row (Binacional 2-0 Alianza Atl., 06/08/2025)
# The date-like text must be forced to text so Excel doesn't convert it
# to a serial date; the leading apostrophe forces text entry, and then
# resetting the style back to "Normal" clears the quote-prefix style
# that Excel applies automatically so the cell keeps the default style.
$ws.Range("A22").Value = "'06/08/2025"
$ws.Range("A22").Style = "Normal"

$ws.Range("B22").Value = "Binacional"
$ws.Range("C22").Value = 2
$ws.Range("D22").Value = 0
$ws.Range("E22").Value = "Alianza Atl."
$ws.Range("F22").Value = "W"
$ws.Range("G22").Value = 2
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 1.5
$ws.Range("L22").Value = 0.51
$ws.Range("M22").Value = 16
$ws.Range("N22").Value = 10
$ws.Range("O22").Value = 7
$ws.Range("P22").Value = 1
